$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at position 200, shifting existing rows 200-286 down to 207-293
$ws.Rows("200:206").Insert()

# Row 200: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A200').Value = 2
$ws.Range('B200').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C200').Value = 'Coquimbo'
$ws.Range('D200').Value = 44524
$ws.Range('E200').Value = 4
$ws.Range('F200').Value = 'Fruta'
$ws.Range('G200').Value = 100106
$ws.Range('H200').Value = 'Oleaginosos'
$ws.Range('I200').Value = 100106002
$ws.Range('J200').Value = 'Palta'
$ws.Range('K200').Value = 'Edranol'
$ws.Range('L200').Value = 'Especial'
$ws.Range('M200').Value = 300
$ws.Range('N200').Value = 1900
$ws.Range('O200').Value = 2000
$ws.Range('P200').Value = 1950
$ws.Range('Q200').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R200').Value = 'Provincia de Limarí'
$ws.Range('S200').Value = 1950
$ws.Range('T200').Value = 1

# Row 201: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A201').Value = 2
$ws.Range('B201').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C201').Value = 'Coquimbo'
$ws.Range('D201').Value = 44524
$ws.Range('E201').Value = 4
$ws.Range('F201').Value = 'Fruta'
$ws.Range('G201').Value = 100106
$ws.Range('H201').Value = 'Oleaginosos'
$ws.Range('I201').Value = 100106002
$ws.Range('J201').Value = 'Palta'
$ws.Range('K201').Value = 'Edranol'
$ws.Range('L201').Value = 'Primera'
$ws.Range('M201').Value = 310
$ws.Range('N201').Value = 1700
$ws.Range('O201').Value = 1800
$ws.Range('P201').Value = 1755
$ws.Range('Q201').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R201').Value = 'Provincia de Limarí'
$ws.Range('S201').Value = 1755
$ws.Range('T201').Value = 1

# Row 202: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A202').Value = 2
$ws.Range('B202').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C202').Value = 'Coquimbo'
$ws.Range('D202').Value = 44524
$ws.Range('E202').Value = 4
$ws.Range('F202').Value = 'Fruta'
$ws.Range('G202').Value = 100106
$ws.Range('H202').Value = 'Oleaginosos'
$ws.Range('I202').Value = 100106002
$ws.Range('J202').Value = 'Palta'
$ws.Range('K202').Value = 'Edranol'
$ws.Range('L202').Value = 'Segunda'
$ws.Range('M202').Value = 200
$ws.Range('N202').Value = 1300
$ws.Range('O202').Value = 1400
$ws.Range('P202').Value = 1350
$ws.Range('Q202').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R202').Value = 'Provincia de Limarí'
$ws.Range('S202').Value = 1350
$ws.Range('T202').Value = 1

# Row 203: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A203').Value = 2
$ws.Range('B203').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C203').Value = 'Coquimbo'
$ws.Range('D203').Value = 44524
$ws.Range('E203').Value = 4
$ws.Range('F203').Value = 'Fruta'
$ws.Range('G203').Value = 100106
$ws.Range('H203').Value = 'Oleaginosos'
$ws.Range('I203').Value = 100106002
$ws.Range('J203').Value = 'Palta'
$ws.Range('K203').Value = 'Hass'
$ws.Range('L203').Value = 'Especial'
$ws.Range('M203').Value = 300
$ws.Range('N203').Value = 2300
$ws.Range('O203').Value = 2400
$ws.Range('P203').Value = 2350
$ws.Range('Q203').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R203').Value = 'Provincia de Limarí'
$ws.Range('S203').Value = 2350
$ws.Range('T203').Value = 1

# Row 204: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A204').Value = 2
$ws.Range('B204').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C204').Value = 'Coquimbo'
$ws.Range('D204').Value = 44524
$ws.Range('E204').Value = 4
$ws.Range('F204').Value = 'Fruta'
$ws.Range('G204').Value = 100106
$ws.Range('H204').Value = 'Oleaginosos'
$ws.Range('I204').Value = 100106002
$ws.Range('J204').Value = 'Palta'
$ws.Range('K204').Value = 'Hass'
$ws.Range('L204').Value = 'Primera'
$ws.Range('M204').Value = 300
$ws.Range('N204').Value = 2000
$ws.Range('O204').Value = 2100
$ws.Range('P204').Value = 2050
$ws.Range('Q204').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R204').Value = 'Provincia de Limarí'
$ws.Range('S204').Value = 2050
$ws.Range('T204').Value = 1

# Row 205: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A205').Value = 2
$ws.Range('B205').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C205').Value = 'Coquimbo'
$ws.Range('D205').Value = 44524
$ws.Range('E205').Value = 4
$ws.Range('F205').Value = 'Fruta'
$ws.Range('G205').Value = 100106
$ws.Range('H205').Value = 'Oleaginosos'
$ws.Range('I205').Value = 100106002
$ws.Range('J205').Value = 'Palta'
$ws.Range('K205').Value = 'Hass'
$ws.Range('L205').Value = 'Segunda'
$ws.Range('M205').Value = 300
$ws.Range('N205').Value = 1700
$ws.Range('O205').Value = 1800
$ws.Range('P205').Value = 1750
$ws.Range('Q205').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R205').Value = 'Provincia de Limarí'
$ws.Range('S205').Value = 1750
$ws.Range('T205').Value = 1

# Row 206: new weekly entry dated 2021-11-24 (serial 44524)
$ws.Range('A206').Value = 2
$ws.Range('B206').Value = 'Comercializadora del Agro de Limarí'
$ws.Range('C206').Value = 'Coquimbo'
$ws.Range('D206').Value = 44524
$ws.Range('E206').Value = 4
$ws.Range('F206').Value = 'Fruta'
$ws.Range('G206').Value = 100106
$ws.Range('H206').Value = 'Oleaginosos'
$ws.Range('I206').Value = 100106002
$ws.Range('J206').Value = 'Palta'
$ws.Range('K206').Value = 'Hass'
$ws.Range('L206').Value = 'Tercera'
$ws.Range('M206').Value = 240
$ws.Range('N206').Value = 1300
$ws.Range('O206').Value = 1400
$ws.Range('P206').Value = 1350
$ws.Range('Q206').Value = '$/kilo (en caja de 17 kilos)'
$ws.Range('R206').Value = 'Provincia de Limarí'
$ws.Range('S206').Value = 1350
$ws.Range('T206').Value = 1
